# Update cryptos list (prices / volume / coin swaps) to match latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "87.757.59"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +6.95%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.308.98"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +3.69%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.02"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.50%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "631.76"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +1.30%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.399"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +36.83%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.649"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +10.79%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.00"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.02%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "3.303.25"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +3.56%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.597"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.92%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +3.00%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +7.25%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.90"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +9.08%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.918.11"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +3.55%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.32"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.16%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "87.532.31"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +6.91%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.304.53"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +3.26%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.37"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.95%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.04"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -6.70%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.29"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +3.21%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "440.74"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.03%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +5.80%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.29"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.43%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.45"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +11.13%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "5.28"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -2.14%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.468.40"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +3.10%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "77.56"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.16%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0000133"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +7.49%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.43%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.185"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +27.78%  "
$ws.Range("B32").NumberFormat = "@"
$ws.Range("B32").Value = "Binance-PegBSC-USD"
$ws.Range("C32").NumberFormat = "@"
$ws.Range("C32").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.01"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.88%  "
$ws.Range("B33").NumberFormat = "@"
$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "9.10"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.04%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "558.41"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -5.57%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -3.29%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.83%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "7.06"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +14.58%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.141"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -9.27%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "22.94"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.48%  "
$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = "FirstDigitalUSD"
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.01%  "
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "WhiteBITCoin"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "21.78"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +4.64%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.68%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.12%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.01"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.90%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.10%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "154.82"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -3.77%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "182.87"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -2.91%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +2.09%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "45.33"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +1.37%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.32"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +2.03%  "
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = "ARBITRUM"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.641"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.45%  "
